$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New job-book rows captured in the "integrated status" sync (rows 63-65).
# Text columns are entered with a leading apostrophe so values that look
# numeric/date-like (incoming numbers, m/d/yyyy strings) are stored as text,
# matching how the rest of the sheet already stores them. The style is then
# reset to "Normal" so no stray number-format is left behind on the cell.
$rows = @(
    @{ Row = 63; A = 71315; B = "6DH13500TSK";  C = "ET  00H01487S"; D = "MULTI"; E = "?"; F = "GOOD"; G = "82934"; H = "NO"; I = ""; J = "ravi"; K = "9/12/2022"; L = $false; M = "N/A"; N = "N/A"; O = "NO" },
    @{ Row = 64; A = 71316; B = "6DP3R35METSK"; C = "21B62598R";      D = "MULTI"; E = "?"; F = "GOOD"; G = "82926"; H = "NO"; I = ""; J = "ravi"; K = "9/12/2022"; L = $false; M = "N/A"; N = "N/A"; O = "NO" },
    @{ Row = 65; A = 71317; B = "DXX36";         C = "3K9X067755";    D = "460";   E = "?"; F = "?";    G = "82867"; H = "NO"; I = ""; J = "ravi"; K = "9/12/2022"; L = $false; M = "N/A"; N = "N/A"; O = "NO" }
)

$textCols = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 13, 14, 15)

foreach ($r in $rows) {
    $row = $r.Row

    # A: jobNumber - numeric
    $ws.Cells.Item($row, 1).Value = $r.A

    # B..K, M..O: text columns (forced via leading apostrophe)
    $ws.Cells.Item($row, 2).Value = "'" + $r.B
    $ws.Cells.Item($row, 3).Value = "'" + $r.C
    $ws.Cells.Item($row, 4).Value = "'" + $r.D
    $ws.Cells.Item($row, 5).Value = "'" + $r.E
    $ws.Cells.Item($row, 6).Value = "'" + $r.F
    $ws.Cells.Item($row, 7).Value = "'" + $r.G
    $ws.Cells.Item($row, 8).Value = "'" + $r.H
    $ws.Cells.Item($row, 9).Value = "'" + $r.I
    $ws.Cells.Item($row, 10).Value = "'" + $r.J
    $ws.Cells.Item($row, 11).Value = "'" + $r.K

    # L: _isDeleted - boolean
    $ws.Cells.Item($row, 12).Value = $r.L

    $ws.Cells.Item($row, 13).Value = "'" + $r.M
    $ws.Cells.Item($row, 14).Value = "'" + $r.N
    $ws.Cells.Item($row, 15).Value = "'" + $r.O

    # Drop the "text" number-format stamp picked up from the apostrophe entry
    # so the cells stay on the default/Normal style, like the rest of the sheet.
    foreach ($c in $textCols) {
        $ws.Cells.Item($row, $c).Style = "Normal"
    }
}
